$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "6-5="
$t.Cell(1,2).Range.Text = "49-17="
$t.Cell(1,3).Range.Text = "73+4="
$t.Cell(1,4).Range.Text = "29+31="
$t.Cell(1,5).Range.Text = "50+12="
$t.Cell(2,1).Range.Text = "33-27="
$t.Cell(2,2).Range.Text = "44-31="
$t.Cell(2,3).Range.Text = "96-75="
$t.Cell(2,4).Range.Text = "37+21="
$t.Cell(2,5).Range.Text = "6+33="
$t.Cell(3,1).Range.Text = "52-31="
$t.Cell(3,2).Range.Text = "92-68="
$t.Cell(3,3).Range.Text = "52-14="
$t.Cell(3,4).Range.Text = "29+60="
$t.Cell(3,5).Range.Text = "65+20="
$t.Cell(4,1).Range.Text = "6+5="
$t.Cell(4,2).Range.Text = "77-2="
$t.Cell(4,3).Range.Text = "45+42="
$t.Cell(4,4).Range.Text = "24+54="
$t.Cell(4,5).Range.Text = "4+42="
$t.Cell(5,1).Range.Text = "77-55="
$t.Cell(5,2).Range.Text = "37-26="
$t.Cell(5,3).Range.Text = "82-76="
$t.Cell(5,4).Range.Text = "70-26="
$t.Cell(5,5).Range.Text = "75-30="
$t.Cell(6,1).Range.Text = "82-38="
$t.Cell(6,2).Range.Text = "78-74="
$t.Cell(6,3).Range.Text = "50+0="
$t.Cell(6,4).Range.Text = "54-41="
$t.Cell(6,5).Range.Text = "68-30="
$t.Cell(7,1).Range.Text = "84-23="
$t.Cell(7,2).Range.Text = "22+8="
$t.Cell(7,3).Range.Text = "91-55="
$t.Cell(7,4).Range.Text = "59+10="
$t.Cell(7,5).Range.Text = "93-80="
$t.Cell(8,1).Range.Text = "38+10="
$t.Cell(8,2).Range.Text = "13+46="
$t.Cell(8,3).Range.Text = "19-2="
$t.Cell(8,4).Range.Text = "37-24="
$t.Cell(8,5).Range.Text = "91-78="
$t.Cell(9,1).Range.Text = "72+14="
$t.Cell(9,2).Range.Text = "32+15="
$t.Cell(9,3).Range.Text = "3+63="
$t.Cell(9,4).Range.Text = "75+8="
$t.Cell(9,5).Range.Text = "43+1="
$t.Cell(10,1).Range.Text = "12+18="
$t.Cell(10,2).Range.Text = "65-56="
$t.Cell(10,3).Range.Text = "68+10="
$t.Cell(10,4).Range.Text = "76+8="
$t.Cell(10,5).Range.Text = "85-7="
$t.Cell(11,1).Range.Text = "35-0="
$t.Cell(11,2).Range.Text = "75-49="
$t.Cell(11,3).Range.Text = "87-9="
$t.Cell(11,4).Range.Text = "42-13="
$t.Cell(11,5).Range.Text = "18+15="
$t.Cell(12,1).Range.Text = "5+33="
$t.Cell(12,2).Range.Text = "93-20="
$t.Cell(12,3).Range.Text = "66-60="
$t.Cell(12,4).Range.Text = "80-6="
$t.Cell(12,5).Range.Text = "49+1="
$t.Cell(13,1).Range.Text = "33-21="
$t.Cell(13,2).Range.Text = "55+4="
$t.Cell(13,3).Range.Text = "89-14="
$t.Cell(13,4).Range.Text = "95-13="
$t.Cell(13,5).Range.Text = "66+25="
$t.Cell(14,1).Range.Text = "13+6="
$t.Cell(14,2).Range.Text = "9+29="
$t.Cell(14,3).Range.Text = "67-55="
$t.Cell(14,4).Range.Text = "59+34="
$t.Cell(14,5).Range.Text = "61+9="
$t.Cell(15,1).Range.Text = "55+33="
$t.Cell(15,2).Range.Text = "16+80="
$t.Cell(15,3).Range.Text = "39+48="
$t.Cell(15,4).Range.Text = "5+63="
$t.Cell(15,5).Range.Text = "75-46="
$t.Cell(16,1).Range.Text = "71+19="
$t.Cell(16,2).Range.Text = "9+84="
$t.Cell(16,3).Range.Text = "85+6="
$t.Cell(16,4).Range.Text = "29+31="
$t.Cell(16,5).Range.Text = "44-25="
$t.Cell(17,1).Range.Text = "31+52="
$t.Cell(17,2).Range.Text = "18+46="
$t.Cell(17,3).Range.Text = "78-3="
$t.Cell(17,4).Range.Text = "9+46="
$t.Cell(17,5).Range.Text = "92-66="
$t.Cell(18,1).Range.Text = "99+0="
$t.Cell(18,2).Range.Text = "15+65="
$t.Cell(18,3).Range.Text = "38+8="
$t.Cell(18,4).Range.Text = "1+6="
$t.Cell(18,5).Range.Text = "37+4="
$t.Cell(19,1).Range.Text = "12+56="
$t.Cell(19,2).Range.Text = "31+55="
$t.Cell(19,3).Range.Text = "46-38="
$t.Cell(19,4).Range.Text = "65-3="
$t.Cell(19,5).Range.Text = "83-62="
$t.Cell(20,1).Range.Text = "97-32="
$t.Cell(20,2).Range.Text = "21-4="
$t.Cell(20,3).Range.Text = "73-28="
$t.Cell(20,4).Range.Text = "8-7="
$t.Cell(20,5).Range.Text = "52-30="
